# take_aways.xlsx - add a new Codeforces problem entry (row 11) to the
# take-away notes table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (B11:F11) - mirrors the existing table layout:
# B=problem link, C=difficulty, D=status, E=first encounter, F=take aways
$ws.Range("B11").Value = "https://codeforces.com/problemset/problem/2002/C"
$ws.Range("C11").Value = "cf 1200"
$ws.Range("D11").Value = "done"
$ws.Range("E11").Value = "looked very hard in the start but had sweet maths solution to it"
$ws.Range("F11").Value = "do not convert numbers into decimal will loose precision so try to keep it as whole numbers in any comparison if possible change equation of conversion from division to multiplication"

# Match the wrap-text/top-aligned look of the rest of the table and size
# the row the way the other multi-line rows are sized.
$ws.Range("B11:F11").WrapText = $true
$ws.Range("B11:F11").VerticalAlignment = -4160
$ws.Rows.Item(11).RowHeight = 72

# Update the on-screen selection / scroll position to reflect where the
# author ended up after entering the new row.
[void]$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "take_aways.xlsx: appended row 11 (codeforces 2002C take-away)"
